$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (the "Undecided" / missing-value row) with the new values
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 55
